# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price refresh values to the Aegis_Profits workbook
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 2028.7778
$ws.Range("J62").Value = 1909.6666
$ws.Range("K62").Value = 2028.7778
$ws.Range("L62").Value = 1909.6666
$ws.Range("M62").Value = -1404.7778
$ws.Range("N62").Value = -3157.6666
# Row 65
$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 2028.7778
$ws.Range("J65").Value = 1909.6666
$ws.Range("K65").Value = 10143.889
$ws.Range("L65").Value = 9548.333000000001
$ws.Range("M65").Value = -7023.889000000001
$ws.Range("N65").Value = -15788.333
# Row 86
$ws.Range("H86").Value = 6731.852
$ws.Range("I86").Value = 6586.875
$ws.Range("J86").Value = 6942.727
$ws.Range("K86").Value = 6586.875
$ws.Range("L86").Value = 6942.727
$ws.Range("M86").Value = -5463.875
$ws.Range("N86").Value = -9188.726999999999
# Row 89
$ws.Range("H89").Value = 6731.852
$ws.Range("I89").Value = 6586.875
$ws.Range("J89").Value = 6942.727
$ws.Range("K89").Value = 32934.375
$ws.Range("L89").Value = 34713.635
$ws.Range("M89").Value = -27318.375
$ws.Range("N89").Value = -45945.635
# Row 111
$ws.Range("H111").Value = 4005538.5
$ws.Range("I111").Value = 10254.272
$ws.Range("J111").Value = 7144690.5
$ws.Range("K111").Value = 30762.816
$ws.Range("L111").Value = 21434071.5
$ws.Range("M111").Value = -27695.816
$ws.Range("N111").Value = -21440205.5
# Row 116
$ws.Range("H116").Value = 1792.7273
$ws.Range("I116").Value = 1297.1111
$ws.Range("K116").Value = 1297.1111
$ws.Range("M116").Value = 2144.8889
# Row 129
$ws.Range("H129").Value = 3050.7446
$ws.Range("J129").Value = 987.45
$ws.Range("L129").Value = 2962.35
$ws.Range("N129").Value = -12962.35
# Row 135
$ws.Range("H135").Value = 489.8
$ws.Range("I135").Value = 489.8
$ws.Range("K135").Value = 4408.2
$ws.Range("M135").Value = -1873.2

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 29119.584
$ws.Range("I2").Value = 1034.7084
$ws.Range("J2").Value = 85289.336
$ws.Range("K2").Value = 1034.7084
$ws.Range("L2").Value = 85289.336
$ws.Range("M2").Value = -921.7084
$ws.Range("N2").Value = -85515.336
# Row 32
$ws.Range("H32").Value = 24624.607
$ws.Range("I32").Value = 4069.242
$ws.Range("K32").Value = 4069.242
$ws.Range("M32").Value = -3782.242
# Row 112
$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22954
# Row 116
$ws.Range("H116").Value = 29119.584
$ws.Range("I116").Value = 1034.7084
$ws.Range("J116").Value = 85289.336
$ws.Range("K116").Value = 1034.7084
$ws.Range("L116").Value = 85289.336
$ws.Range("M116").Value = 1259.2916
$ws.Range("N116").Value = -89877.336

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 29119.584
$ws.Range("I3").Value = 1034.7084
$ws.Range("J3").Value = 85289.336
$ws.Range("K3").Value = 1034.7084
$ws.Range("L3").Value = 85289.336
$ws.Range("M3").Value = -920.7084
$ws.Range("N3").Value = -85517.336

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 1912.5
$ws.Range("J17").Value = 1650
$ws.Range("L17").Value = 4950
$ws.Range("N17").Value = -5288
# Row 107
$ws.Range("H107").Value = 253504.64
$ws.Range("I107").Value = 398.66666
$ws.Range("J107").Value = 633163.5600000001
$ws.Range("K107").Value = 1195.99998
$ws.Range("L107").Value = 1899490.68
$ws.Range("M107").Value = 724.0000199999999
$ws.Range("N107").Value = -1903330.68
# Row 131
$ws.Range("H131").Value = 1499.0204
$ws.Range("J131").Value = 1465.8223
$ws.Range("L131").Value = 4397.4669
$ws.Range("N131").Value = -14477.4669

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 4348333.5
$ws.Range("I11").Value = 5018000
$ws.Range("J11").Value = 1000000
$ws.Range("K11").Value = 5018000
$ws.Range("L11").Value = 1000000
$ws.Range("M11").Value = -5017861
$ws.Range("N11").Value = -1000278
# Row 43
$ws.Range("H43").Value = 4219.6924
$ws.Range("J43").Value = 5228.4443
$ws.Range("L43").Value = 5228.4443
$ws.Range("N43").Value = -5530.4443
# Row 46
$ws.Range("H46").Value = 8649.666999999999
$ws.Range("I46").Value = 7474.5
$ws.Range("J46").Value = 11000
$ws.Range("K46").Value = 7474.5
$ws.Range("L46").Value = 11000
$ws.Range("M46").Value = -7318.5
$ws.Range("N46").Value = -11312
# Row 57
$ws.Range("H57").Value = 15250
$ws.Range("J57").Value = 15250
$ws.Range("L57").Value = 15250
$ws.Range("N57").Value = -16890
# Row 70
$ws.Range("H70").Value = 82438.03999999999
$ws.Range("I70").Value = 130124.625
$ws.Range("J70").Value = 6139.5
$ws.Range("K70").Value = 130124.625
$ws.Range("L70").Value = 6139.5
$ws.Range("M70").Value = -129854.625
$ws.Range("N70").Value = -6679.5
# Row 73
$ws.Range("H73").Value = 82438.03999999999
$ws.Range("I73").Value = 130124.625
$ws.Range("J73").Value = 6139.5
$ws.Range("K73").Value = 130124.625
$ws.Range("L73").Value = 6139.5
$ws.Range("M73").Value = -129188.625
$ws.Range("N73").Value = -8011.5
# Row 80
$ws.Range("H80").Value = 3920.8
$ws.Range("I80").Value = 5201.6665
$ws.Range("K80").Value = 5201.6665
$ws.Range("M80").Value = -4203.6665
# Row 83
$ws.Range("H83").Value = 3920.8
$ws.Range("I83").Value = 5201.6665
$ws.Range("K83").Value = 26008.3325
$ws.Range("M83").Value = -21016.3325
# Row 97
$ws.Range("H97").Value = 58825950
$ws.Range("I97").Value = 83336040
$ws.Range("J97").Value = 1742
$ws.Range("K97").Value = 83336040
$ws.Range("L97").Value = 1742
$ws.Range("M97").Value = -83335544
$ws.Range("N97").Value = -2734
# Row 102
$ws.Range("H102").Value = 431966.56
$ws.Range("I102").Value = 5177
$ws.Range("K102").Value = 5177
$ws.Range("M102").Value = -3555
# Row 111
$ws.Range("H111").Value = 36646.5
$ws.Range("J111").Value = 36646.5
$ws.Range("L111").Value = 36646.5
$ws.Range("N111").Value = -42780.5

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 168894.67
$ws.Range("I40").Value = 251472
$ws.Range("K40").Value = 251472
$ws.Range("M40").Value = -251336
# Row 100
$ws.Range("H100").Value = 1489.7333
$ws.Range("I100").Value = 1433.0834
$ws.Range("J100").Value = 1716.3334
$ws.Range("K100").Value = 1433.0834
$ws.Range("L100").Value = 1716.3334
$ws.Range("M100").Value = -892.0834
$ws.Range("N100").Value = -2798.3334
# Row 110
$ws.Range("H110").Value = 30800
$ws.Range("J110").Value = 30800
$ws.Range("L110").Value = 30800
$ws.Range("N110").Value = -38980
# Row 127
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920
# Row 132
$ws.Range("H132").Value = 8175.5
$ws.Range("I132").Value = 12433.777
$ws.Range("J132").Value = 2700.5715
$ws.Range("K132").Value = 37301.331
$ws.Range("L132").Value = 8101.7145
$ws.Range("M132").Value = -34771.331
$ws.Range("N132").Value = -13161.7145

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2295.182
$ws.Range("I126").Value = 2467.75
$ws.Range("J126").Value = 1835
$ws.Range("K126").Value = 7403.25
$ws.Range("L126").Value = 5505
$ws.Range("M126").Value = -4933.25
$ws.Range("N126").Value = -10445
# Row 128
$ws.Range("H128").Value = 48997.6
$ws.Range("J128").Value = 48997.6
$ws.Range("L128").Value = 48997.6
$ws.Range("N128").Value = -58957.6
# Row 132
$ws.Range("H132").Value = 6630.7144
$ws.Range("I132").Value = 6756.231
$ws.Range("K132").Value = 20268.693
$ws.Range("M132").Value = -17738.693
# Row 133
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120
